$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.528.29"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "3.103.10"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "386.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.539"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0854"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").Value = "3.590.83"
$ws.Range("E13").Value = "  +2.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "3.098.72"
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.992"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.43%  "
$ws.Range("D19").Value = "51.572.25"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.90%  "
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").Value = "0.0₃0963"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.70%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.166"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0479"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.14%  "
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.291"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "129.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.116"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.68%  "
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.36%  "
$ws.Range("E47").Value = "  +5.11%  "
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").Value = "2.064.10"
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.938"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +18.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0329"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.79%  "
